$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ranking scores (B2:B4) must stay as TEXT (they were already stored as
# inline-string text, e.g. "0.091"), so force Text format before writing
# the numeric-looking strings -- otherwise Excel auto-converts them to numbers.
$ws.Range("B2:B4").NumberFormat = "@"

# Row 2 (was Noticia N10 -> now Noticia N09 content)
$ws.Range("B2").Value = '0.109'
$ws.Range("C2").Value = 'Noticia N° 09.txt'
$ws.Range("D2").Value = 'titulo: aport dol soja, banc central compr sol usd millon merc resumen: entid monetari encaden rued oper compr netas, marg achic prim dia agreg tip cambi especial export sojer contenido: jun, sesion usd millon segment cont spot, liquid program increment export econom regionales, dolar, aport usd millones, volum aport vent soj derivados, goz estand cambiari seman pasada. bcra conclu intervencion cambiari dia sald net comprador sol millon dolares. oper dol agro registr tercer etap ingres usd millon abril parte. asimismo, empez balanc bcra intervencion cambiari negat usd millones, mientr juni sostien im sald comprador millon dolares. bcra aceler resguard escas reservas, med afect provinci municipios, moment ministeri econom plane nuev canj voluntari bon pes objet despej vencimientos. ministr economia, sergi massa, viaj washington cerc juni procur cerr readecu acuerd fond monetari internacional fmi, permit argentin obten adelant desembols nuev met cumplir, lueg impact sequ export sector agro, “las negoci argentin funcionari tecnic fond avanz hac casi dos mes traves reunion virtuales, defin modific viej acuerdo, virtual suspend lueg cumpl met reserv fiscal prim trimestr ″, coment expert research traders. objet econom fond adelante, menos, part desembols compromet fin año usd millones, ayud reforz reserv año export sector agro caer cerc usd millones, acuerd proyeccion bols comerci rosari bcr. tal desembols dud cuant podr utiliz interven merc cambiario. fond acept bcra interveng eventual situacion stress. todav resuelt mont fmi consent intervenciones. reserv internacional brut banc central crec seman pas usd millon finaliz millon dolares. inform anker latinoamer subray cuant reserv liqu bcra, “su disponibil hoy dad gran med encaj cuent bancari moned deposit bcra -usd millones-. estabil deposit moned extranjer crucial sosten capac intervencion bcra”. lueg divers gestion ministr mass incentiv pag import divis china, expand uso yuan comerci exterior. ener may oper equivalent usd millones. ademas, dos empres fabric electron confirm pag compromis deud total usd millon moneda, alivi arcas bcra super usd millon prim trimestr segu leyendo: urls imagenes:'

# Row 3 (was Noticia N05 -> now Noticia N08 content)
$ws.Range("B3").Value = '0.086'
$ws.Range("C3").Value = 'Noticia N° 08.txt'
$ws.Range("D3").Value = 'titulo: opin econom crecimient creacion emple tras polem lousteau cerruti resumen: entrev infobae, senador precandidat jef gobiern asegur econom argentin cre trabaj hac años crece. portavoz presidencial cruzo; dic anal contenido: jun, ximen cas portavoz presidencial gabriel cerruti precandidat jef gobiern porteñ junt cambi martin lousteau cruz traves red social distint vision crecimient creacion emple argentin ultim decadas. entrev infobae, econom lousteau asegur econom argentin “est trabada”, cre trabaj hac años crece. lueg publicacion, cerruti respond seri tuits destac van mes crecimient consecut trabaj registrado, period crecimient prolong ultim años. “en mes incorpor trabaj trabaj empres privadas”, dij funcionaria. “ant afirm pais crec acerc datos: pbi años atras ultim dato, mil millon dolares, dat crec mil millones. si crecimiento”, detall cerruti. lousteau desmint cifras: “est numer implic ingres per capit argentin crec aument tas promedi año. argentin puest ranking pais mayor crecimient periodo, junt india, malasia, singapur, core china, otros”. ¿cual vision econom creacion emple argentina? “mas alla numer exactos, punt central emple cre ultim años casi sector public sector priv formal. pued hab crec practic nada. mayor part emple gener sector public traves monotribut emple formal. clar econom dinam cre emple calidad. cre emple precari sector publico, segur encubiert desempleo”, resalt econom miguel kiguel, director econviews integr centr estabil financ cef. juan luis bour, director econom jef fiel, desmint numer argument present cerruti: ultim dat oficial public febrer muestr emple formal priv epf crec puest años, febrer febrer apen puest año. mism periodo, mism estadist ministeri trabaj revel emple public crec puestos. cuent cerruti comp ultim dat peor numer pandemia: juni perd puestos. “la recuper puest da ‘normalizacion’ econom lueg profund caida. tipic error algui sab leer estadisticas”, explico. problem central emple casi tod creacion net puest informales, cuentaprop emple publicos. econom crec %, nuev emple informales, cuentaprop emple publicos. argentin segun dat indec cuent gener ingres casi ocup contrat laboral informal cuentapropistas. expuest fuert oscil emple -son ajust pandemia- perd termin ingres real ajust inflacion. asalari priv formal numer bastant establ crec sujet violent oscil ingres real volatil inflacion. grup “protegido” asalari emple publicos, crec form establ tas superior poblacion cuy ingres volv crec termin reales. lado, bour cuestion argument utiliz lousteau. “los comentari lousteau emple claros, ello contribu discusion, aunqu entiend ambient formulados. problem emple –si refier emple priv calidad, decir, asalari formal- sol tem crecimient economico”, dijo. econom fiel destac men dos variabl ausent diagnost ex ministr economia. “la primer tod norm regulatori anti-empleo, nadi contrat aun crecimient sab dram sobrevien si cae demanda; deb cerr fabrica, pued reestructurar”, aseguro. “la segund si econom estabiliz inflacion anual, mensual, empres ahor pod licu cost inflacion deb aprend conviv escenari demand much flexibil convertibil problema. mism fisc deb aprend viv pod licu gast inflacion. mund diferent nadi parec plant requier introduc flexibil pod reestructurar”, conclu bour segun indec, emple informal, cobertur social descuent jubilatorio, ubic nivel record total. result lleg despues años sol mit emple cre correspond puest trabaj formales, dentr sector privado, segun estudi realiz institut desarroll social argentin idesa. acuerd calcul institut privado, ultim dos dec cre millon puest trabajo, cual millon gener sector publico, millon correspond emple informal restant millon cre sector privado. dat muestr argentin debil estructural gener emple calidad. primer consecuent años logr reduc gran cantid adult edad activ trabaj inact desempleados. deriv casi mit gent edad trabaj permanec merc trabajo”, indic ides reporte. segu leyendo: urls imagenes:'

# Row 4 (was Noticia N09 -> now Noticia N06 content)
$ws.Range("B4").Value = '0.074'
$ws.Range("C4").Value = 'Noticia N° 06.txt'
$ws.Range("D4").Value = 'titulo: inflacion may argentin super venezuela, segun estim priv resumen: abril registr resultado; vari preci ultim año ambos pais esper diciembr contenido: jun, martin kanenguis inflacion argentin may super venezuela, vez mas, segun estim priv ambos paises. dat observatori venezolan finanz ovf, sig anal deb falt rigor estadist oficiales, marc sub preci mes pasado, desaceler ultim meses. argentina, segun mayor estimaciones, dat mes pas rond %, mientr gobiern rez ubiqu levement debaj cifra. abril, argentin super venezuela: registr %, frent inform ovf inform banc central venezuela, dab conoc inform octubr año pasado. ultim mes mayo-mayo, segun ovf, inflacion pais gobern nicolas madur lleg %, baj respect registr abril. ademas, ener sub preci acumul %. tendenci inflacion ultim mes divergentes: marz ovf registr venezuel %, abril may %; tanto, indec registr argentin marzo, abril estim mes qued torn interanual; conoc dat organ lider marc lavagna. segun ovf, “en medi econom evident sign contraccion, tas inflacion da tregua. asi, may tas inflacion mensual alcanz %, acumul anualiz %. respect abril, tas mensual inflacion triplic aunqu tas doc mes desacelero”. “est comport inflacion ocurr entorn sign signific caid nivel activ indiscutibles, ostens menor vent comerci produccion industrial contraj prim trimestr ″, indico. ademas, “la debil demand agreg notori deb salari pension sector public pronunci rezag respect inflacion ejecu gast part gobiern baja”. “al compar cifr inflacion may respect aument tip cambi mism mes, clar estrech relacion ambas variables, obstante, menciono, debil demanda, explic polit compresion salarial ven aplic gobierno”, subray ovf. cuant principal component conform indic nacional preci consumidor, destac “los aument experiment rubr esparc %, vest calz %, equip hog alquil viviend %. aliment alzas modestas, may increment %”. “tod ello pon manifiesto, polit econom aplic sid ineficaz conten alza precios, aun retraccion induc demand agregada”, conclu organ independient regim autoritari maduro. abril, banc central venezuel inform aument cuatr meses, lueg difund dat medi año. tanto, mayor relev preci consumidor consultor argentin anticip ipc torn %. c&t indic relev preci minor region gba “present alza mensual %, super larg variacion abril may año pasado. asi, variacion doc mes trep %, mayor agost ″. “el rubr mayor increment viviend %, reflej sub gas electr principalmente”, aclaro. bien servici vari “ocup segund lugar, alza %, explic cigarrill articul tocador”. tanto, “el comport esparc fuertement influ alza dolar financier fin abril, vio reflej turism product electronicos”. vez, “aliment beb crec mes. arranc mes gran impuls lueg moderando. verduras, lacte deriv harin destacaron, igual aliment consum hog llevar”. “en salud destac increment medicamentos, sum nuev ajust prepagas”, indic c&t. “en equip manten hogar, artefact hog rol preponder man alza dolar financier pes liquid mes”, concluyo. parte, ecog inform si bien ultim seman mes pas registr fuert desaceler lug variacion aliment respect seman anterior, termin general inflacion sid ciento. particular, aument preci aliment sid ciento. “si consider ademas evolu aliment consum hog %, inflacion aliment alcanz %”, aclar consultor dirig marin dal poggeto. lado, lcg detall “el indic aliment beb present inflacion mensual promedi ultim seman punt punt mism periodo”. mes sum aument preci servici transporte, prepagas, colegi priv combustibles, valor regulados. equip econom trat control cuestion cambiari inflacion aceler todav motiv res ped fmi aceler devalu tip cambi oficial. segu leyendo: urls imagenes:'

# Column D was widened (re-autofit) to fit the new, longer "Contenido" text
$ws.Columns.Item(4).ColumnWidth = 4742.8
